# Update the legacy GSC export data:
#  - Drop the oldest date row (2025-10-01), which shifts every remaining
#    row's date and URL counts up by one row.
#  - Append a new row for the next day (2025-12-30) at the bottom of the
#    table, with the URL counts reset to 0, matching the existing pattern
#    used for days without recorded traffic yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the 2025-10-01 row (row 2); everything below shifts up by one.
$ws.Rows.Item(2).Delete() | Out-Null

# The table now ends at row 90 (2025-12-29). Append the new day, 2025-12-30.
$lastRow = $ws.UsedRange.Rows.Count + 1

$dateCell = $ws.Cells.Item($lastRow, 1)
# Force text so the date string isn't auto-converted to a date serial value.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-30"

$ws.Cells.Item($lastRow, 2).Value = 0.0
$ws.Cells.Item($lastRow, 3).Value = 0.0
